$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Change 1: the "Kommentar:" paragraph about fridlysta arter (first
# occurrence, under "Ur FSC-standarden") loses its italic run
# formatting and gains a trailing period.
# --------------------------------------------------------------------
$old1 = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen"
$new1 = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen."

$rng1 = $d.Content
$found1 = $rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Text = $new1
    $rng1.Font.Italic = 0
}

# --------------------------------------------------------------------
# Change 2: the "6.4" paragraph loses its trailing space.
# --------------------------------------------------------------------
$old2 = "Certifikatsinnehavaren ska skydda sällsynta arter och hotade arter samt deras livsmiljöer inom skogsbruksenheten. Det ska ske genom avsättningar, andra skyddade områden och genom att skapa konnektivitet och/eller genom andra direkta åtgärder som gynnar dessa arters överlevnad och livskraft. Åtgärderna ska stå i förhållande till brukandets skala, intensitet och risk, samt till sällsynta och hotade arters bevarandestatus och ekologiska krav. Certifikatsinnehavaren ska beakta den geografiska spridningen och ekologiska krav hos sällsynta och hotade arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas. "
$new2 = "Certifikatsinnehavaren ska skydda sällsynta arter och hotade arter samt deras livsmiljöer inom skogsbruksenheten. Det ska ske genom avsättningar, andra skyddade områden och genom att skapa konnektivitet och/eller genom andra direkta åtgärder som gynnar dessa arters överlevnad och livskraft. Åtgärderna ska stå i förhållande till brukandets skala, intensitet och risk, samt till sällsynta och hotade arters bevarandestatus och ekologiska krav. Certifikatsinnehavaren ska beakta den geografiska spridningen och ekologiska krav hos sällsynta och hotade arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas."

$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null


# --------------------------------------------------------------------
# Change 3: remove the "6.4.1 Följande biotoper ..." paragraph and the
# following "b) nyckelbiotoper ..." paragraph, then renumber the
# remaining "6.4.1 Bevarandeåtgärder ..." paragraph to "6.4.3 ".
# --------------------------------------------------------------------
$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Följande biotoper undantas")) {
        $startPara = $p
    }
}
if ($startPara -ne $null) {
    $delStart = $startPara.Range.Start
    $keepPara = $startPara.Next().Next()
    $delEnd = $keepPara.Range.Start
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}

$renumRng = $d.Content
$foundRenum = $renumRng.Find.Execute("6.4.1 Bevarandeåtgärder genomförs", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundRenum) {
    $boldRng = $d.Range($renumRng.Start, $renumRng.Start + 6)
    Write-Output "boldRng text=[$($boldRng.Text)]"
    $boldRng.Text = "6.4.3 "
}

Write-Output "step3 done found=$foundRenum"

# --------------------------------------------------------------------
# Change 4: update the date in the first-page header.
# --------------------------------------------------------------------
$sec = $d.Sections(1)
$hdr = $sec.Headers(2)
$foundDate = $hdr.Range.Find.Execute("2023-10-22", $false, $false, $false, $false, $false, $true, 1, $false, "2023-10-25", 2)
Write-Output "date done found=$foundDate"

Write-Output "step1+2 done"
